$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in rows 15 and 16 with new contribution log entries (Sprint 4, Will McLain)
$ws.Range("B15").Value = 4
$ws.Range("C15").Value = "Will McLain"
$ws.Range("D15").Value = "Lead meeting while Ryan was out of town"

$ws.Range("B16").Value = 4
$ws.Range("C16").Value = "Will McLain"
$ws.Range("D16").Value = "Assigned duties for sprint 4"

# Update the selection to match the new active cell
$ws.Range("D18").Select()
